## 1 GFG problem on LL and Doubly Linked List
## Adds two new rows to the Linked List question tracker worksheet:
##   - "Flattening a Linked List" (GFG / Java)
##   - "Remove duplicates from a sorted doubly linked list" (GFG / Java)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (green fill etc.) of the last existing data row (35)
# down onto the two new rows before filling in their values, so the new
# rows visually match the rest of the "GFG / Java" entries.
$ws.Range("A35:C35").Copy()
$ws.Range("A36:C37").PasteSpecial(-4122)

$ws.Range("A36").Value = "GFG"
$ws.Range("B36").Value = "Flattening a Linked List"
$ws.Range("C36").Value = "Java"

$ws.Range("A37").Value = "GFG"
$ws.Range("B37").Value = "Remove duplicates from a sorted doubly linked list"
$ws.Range("C37").Value = "Java"

# Restore the cursor/selection position left behind in the saved file.
$ws.Activate()
$ws.Range("B41").Select()
